$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric value updates ---

# CannonCore V2 (column C): Max Power Loss, Open Loop RPM, Sensorless ERPM
$ws.Range("C5").Value = 700
$ws.Range("C6").Value = 1500
$ws.Range("C7").Value = 1500

# Hypercore (column D): Max Motor Current corrected from 150 to 120
$ws.Range("D9").Value = 120

# Min Torque (row 11) for CannonCore V1 (B) and CannonCore V2 (C)
$ws.Range("B11").Value = 48
$ws.Range("C11").Value = 53

# --- New hyperlink: "More info" link for CannonCore V2 ---
$ws.Hyperlinks.Add($ws.Range("C13"), "https://floatwheel.co/index.php?route=product/product&path=84_89&product_id=9913")

# Match the look of the other "More info" hyperlink cells in this row (F13/G13/H13/I13)
$ws.Range("C13").Font.Color = $ws.Range("F13").Font.Color
$ws.Range("C13").Font.Size = $ws.Range("F13").Font.Size
